$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: shift attendance dates, add a new date column (F) ---
# Force these to stay plain text (not auto-converted to date serials) and
# land back on the default (unstyled) cell format, matching the source cells.
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "2025-04-10"
$ws.Range("D1").Style = "Normal"

$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "2025-04-23"
$ws.Range("E1").Style = "Normal"

$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = "2025-04-24"
$ws.Range("F1").Style = "Normal"

# --- Update / extend student roster: A=Nombre, B=Apellido, C=Numero de Alumno ---
$ws.Range("A2").Value = "Armando"
$ws.Range("B2").Value = "Díaz"
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = "Cristian"
$ws.Range("B3").Value = "Quintero"
$ws.Range("C3").Value = 2

$ws.Range("A4").Value = "Javier"
$ws.Range("B4").Value = "Miranda"
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = "Nataly"
$ws.Range("B5").Value = "García"
$ws.Range("C5").Value = 4

$ws.Range("A6").Value = "Angelo"
$ws.Range("B6").Value = "Garcia"
$ws.Range("C6").Value = 5

$ws.Range("A7").Value = "Yoselin"
$ws.Range("B7").Value = "Reyes"
$ws.Range("C7").Value = 7

$ws.Range("A8").Value = "Jesus"
$ws.Range("B8").Value = "Martinez"
$ws.Range("C8").Value = 6

$ws.Range("A9").Value = "Miguel"
$ws.Range("B9").Value = "Collin"
$ws.Range("C9").Value = 8

$ws.Range("A10").Value = "Edgar"
$ws.Range("B10").Value = "Miranda"
$ws.Range("C10").Value = 9

$ws.Range("A11").Value = "Rosas"
$ws.Range("B11").Value = "Santiago"
$ws.Range("C11").Value = 10

$ws.Range("A12").Value = "Cristian"
$ws.Range("B12").Value = "Gabriel"
$ws.Range("C12").Value = 11

$ws.Range("A13").Value = "Adrian"
$ws.Range("B13").Value = "Martínez"
$ws.Range("C13").Value = 12

$ws.Range("A14").Value = "Wendy"
$ws.Range("B14").Value = "Santiago"
$ws.Range("C14").Value = 13

$ws.Range("A15").Value = "Alexis"
$ws.Range("B15").Value = "Miranda"
$ws.Range("C15").Value = 14

$ws.Range("A16").Value = "Javier"
$ws.Range("B16").Value = "Cruz"
$ws.Range("C16").Value = 15

$ws.Range("A17").Value = "Zuriel"
$ws.Range("B17").Value = "Fernando"
$ws.Range("C17").Value = 16

$ws.Range("A18").Value = "Leo"
$ws.Range("B18").Value = "Manuel"
$ws.Range("C18").Value = 18

$ws.Range("A19").Value = "Emmanuel"
$ws.Range("B19").Value = "Medina"
$ws.Range("C19").Value = 17

$ws.Range("A20").Value = "Uriel"
$ws.Range("B20").Value = "Camacho"
$ws.Range("C20").Value = 20

$ws.Range("A21").Value = "Emanuel"
$ws.Range("B21").Value = "Cresensiano"
$ws.Range("C21").Value = 19

$ws.Range("A22").Value = "Adair"
$ws.Range("B22").Value = "Antonio"
$ws.Range("C22").Value = 21

$ws.Range("A23").Value = "Alexis"
$ws.Range("B23").Value = "Hernández"
$ws.Range("C23").Value = 22

$ws.Range("A24").Value = "Emilio"
$ws.Range("B24").Value = "Galvan"
$ws.Range("C24").Value = 24

$ws.Range("A25").Value = "Josue"
$ws.Range("B25").Value = "Gregorio"
$ws.Range("C25").Value = 23

# --- Attendance (D/E) marker cells: row 5 moves its blank marker from D to E,
#     row 2's extra blank marker in E is dropped, row 6 drops its lone E marker ---
$ws.Range("E2").Clear()
$ws.Range("D5").Clear()
$ws.Range("E6").Clear()

$ws.Range("D2").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New "Presente" attendance column (F), styled like the existing blank
#     attendance markers (green fill) for every student row 2-25 ---
$ws.Range("D2").Copy()
$ws.Range("F2:F25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 6).Value = "Presente"
}
